# Auto-generated edit script: append rows 27-36 to 'Report' sheet, adjust print area/selection/styles
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('Report')

function Apply-CenterStyle($rng, $filled) {
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $false
    if ($filled) { $rng.Interior.Color = 14936315 }
    $rng.Borders.LineStyle = 1
}

function Apply-LeftStyle($rng, $filled, $wrap) {
    $rng.HorizontalAlignment = -4131
    $rng.VerticalAlignment = -4108
    if ($wrap) { $rng.WrapText = $true } else { $rng.WrapText = $false }
    if ($filled) { $rng.Interior.Color = 14936315 }
    $rng.Borders.LineStyle = 1
}

# ---- Row 27 ----
Apply-CenterStyle $ws.Range('A27:L27') $true
Apply-CenterStyle $ws.Range('N27') $true
Apply-CenterStyle $ws.Range('Q27:AB27') $true
Apply-CenterStyle $ws.Range('AD27:AK27') $true
Apply-LeftStyle $ws.Range('M27') $true $false
Apply-LeftStyle $ws.Range('O27') $true $false
Apply-LeftStyle $ws.Range('P27') $true $true
Apply-LeftStyle $ws.Range('AC27') $true $true
$ws.Range('A27').Value = 25
$ws.Range('B27').Value = '服務'
$ws.Range('C27').Value = 2025060733
$ws.Range('F27').Value = 3929
$ws.Range('G27').Value = '蘆洲中山一'
$ws.Range('H27').Value = '新北市蘆洲區'
$ws.Range('Q27').Value = 'THILF03929'
$ws.Range('R27').Value = '新北一'
$ws.Range('S27').Value = '吳宗鴻'
$ws.Range('T27').Value = 1
$ws.Range('U27').Value = '已完工'
$ws.Range('V27').Value = '2025-06-05 11:21:29'
$ws.Range('W27').Value = '2025-06-05 11:00:00'
$ws.Range('X27').Value = '2025-06-05 11:20:00'
$ws.Range('Z27').Value = 0.3
$ws.Range('AB27').Value = '到場處理'
$ws.Range('AC27').Value = 'PMQ2+EDC+STAR'
$ws.Range('AD27').Value = 'O'
$ws.Range('AJ27').Value = 'O'
$ws.Range('AK27').Value = 'O'

# ---- Row 28 ----
Apply-CenterStyle $ws.Range('A28:L28') $false
Apply-CenterStyle $ws.Range('N28') $false
Apply-CenterStyle $ws.Range('Q28:AB28') $false
Apply-CenterStyle $ws.Range('AD28:AK28') $false
Apply-LeftStyle $ws.Range('M28') $false $false
Apply-LeftStyle $ws.Range('O28') $false $false
Apply-LeftStyle $ws.Range('P28') $false $true
Apply-LeftStyle $ws.Range('AC28') $false $true
$ws.Range('A28').Value = 26
$ws.Range('B28').Value = '服務'
$ws.Range('C28').Value = 2025060770
$ws.Range('F28').Value = 3452
$ws.Range('G28').Value = '北縣蘆信店'
$ws.Range('H28').Value = '新北市蘆洲區'
$ws.Range('Q28').Value = 'THILF03452'
$ws.Range('R28').Value = '新北一'
$ws.Range('S28').Value = '吳宗鴻'
$ws.Range('T28').Value = 1
$ws.Range('U28').Value = '已完工'
$ws.Range('V28').Value = '2025-06-05 12:44:19'
$ws.Range('W28').Value = '2025-06-05 12:00:00'
$ws.Range('X28').Value = '2025-06-05 12:43:00'
$ws.Range('Z28').Value = 0.7
$ws.Range('AB28').Value = '到場處理'
$ws.Range('AC28').Value = 'PMQ2+EDC+STAR'
$ws.Range('AD28').Value = 'O'
$ws.Range('AJ28').Value = 'O'
$ws.Range('AK28').Value = 'O'

# ---- Row 29 ----
Apply-CenterStyle $ws.Range('A29:L29') $true
Apply-CenterStyle $ws.Range('N29') $true
Apply-CenterStyle $ws.Range('Q29:AB29') $true
Apply-CenterStyle $ws.Range('AD29:AK29') $true
Apply-LeftStyle $ws.Range('M29') $true $false
Apply-LeftStyle $ws.Range('O29') $true $false
Apply-LeftStyle $ws.Range('P29') $true $true
Apply-LeftStyle $ws.Range('AC29') $true $true
$ws.Range('A29').Value = 27
$ws.Range('B29').Value = '服務'
$ws.Range('C29').Value = 2025060781
$ws.Range('F29').Value = 4316
$ws.Range('G29').Value = '五股工商店'
$ws.Range('H29').Value = '新北市五股區'
$ws.Range('Q29').Value = 'THILF04316'
$ws.Range('R29').Value = '新北一'
$ws.Range('S29').Value = '湯家瑋'
$ws.Range('T29').Value = 1
$ws.Range('U29').Value = '已完工'
$ws.Range('V29').Value = '2025-06-05 13:55:27'
$ws.Range('W29').Value = '2025-06-05 13:30:00'
$ws.Range('X29').Value = '2025-06-05 13:50:00'
$ws.Range('Z29').Value = 0.3
$ws.Range('AB29').Value = '到場處理'
$ws.Range('AC29').Value = 'PMQ2+EDC'
$ws.Range('AD29').Value = 'O'
$ws.Range('AJ29').Value = 'O'
$ws.Range('AK29').Value = 'O'

# ---- Row 30 ----
Apply-CenterStyle $ws.Range('A30:L30') $false
Apply-CenterStyle $ws.Range('N30') $false
Apply-CenterStyle $ws.Range('Q30:AB30') $false
Apply-CenterStyle $ws.Range('AD30:AK30') $false
Apply-LeftStyle $ws.Range('M30') $false $false
Apply-LeftStyle $ws.Range('O30') $false $false
Apply-LeftStyle $ws.Range('P30') $false $true
Apply-LeftStyle $ws.Range('AC30') $false $true
$ws.Range('A30').Value = 28
$ws.Range('B30').Value = '服務'
$ws.Range('C30').Value = 2025060792
$ws.Range('F30').Value = 4218
$ws.Range('G30').Value = '蘆洲湧蓮店'
$ws.Range('H30').Value = '新北市蘆洲區'
$ws.Range('Q30').Value = 'THILF04218'
$ws.Range('R30').Value = '新北一'
$ws.Range('S30').Value = '吳宗鴻'
$ws.Range('T30').Value = 1
$ws.Range('U30').Value = '已完工'
$ws.Range('V30').Value = '2025-06-05 14:18:15'
$ws.Range('W30').Value = '2025-06-05 13:00:00'
$ws.Range('X30').Value = '2025-06-05 13:30:00'
$ws.Range('Z30').Value = 0.5
$ws.Range('AB30').Value = '到場處理'
$ws.Range('AC30').Value = 'PMQ2+EDC+STAR'
$ws.Range('AD30').Value = 'O'
$ws.Range('AJ30').Value = 'O'
$ws.Range('AK30').Value = 'O'

# ---- Row 31 ----
Apply-CenterStyle $ws.Range('A31:L31') $true
Apply-CenterStyle $ws.Range('N31') $true
Apply-CenterStyle $ws.Range('Q31:AB31') $true
Apply-CenterStyle $ws.Range('AD31:AK31') $true
Apply-LeftStyle $ws.Range('M31') $true $false
Apply-LeftStyle $ws.Range('O31') $true $false
Apply-LeftStyle $ws.Range('P31') $true $true
Apply-LeftStyle $ws.Range('AC31') $true $true
$ws.Range('A31').Value = 29
$ws.Range('B31').Value = '服務'
$ws.Range('C31').Value = 2025060795
$ws.Range('F31').Value = 'D349'
$ws.Range('G31').Value = '板橋成都店'
$ws.Range('H31').Value = '新北市板橋區'
$ws.Range('Q31').Value = 'THILF0D349'
$ws.Range('R31').Value = '新北一'
$ws.Range('S31').Value = '狄澤洋'
$ws.Range('T31').Value = 1
$ws.Range('U31').Value = '已完工'
$ws.Range('V31').Value = '2025-06-05 14:20:27'
$ws.Range('W31').Value = '2025-06-05 13:58:00'
$ws.Range('X31').Value = '2025-06-05 14:20:00'
$ws.Range('Z31').Value = 0.4
$ws.Range('AB31').Value = '到場處理'
$ws.Range('AC31').Value = 'PMQ2+EDC'
$ws.Range('AD31').Value = 'O'
$ws.Range('AJ31').Value = 'O'
$ws.Range('AK31').Value = 'O'

# ---- Row 32 ----
Apply-CenterStyle $ws.Range('A32:L32') $false
Apply-CenterStyle $ws.Range('N32') $false
Apply-CenterStyle $ws.Range('Q32:AB32') $false
Apply-CenterStyle $ws.Range('AD32:AK32') $false
Apply-LeftStyle $ws.Range('M32') $false $false
Apply-LeftStyle $ws.Range('O32') $false $false
Apply-LeftStyle $ws.Range('P32') $false $true
Apply-LeftStyle $ws.Range('AC32') $false $true
$ws.Range('A32').Value = 30
$ws.Range('B32').Value = '服務'
$ws.Range('C32').Value = 2025060796
$ws.Range('F32').Value = 2958
$ws.Range('G32').Value = '北縣蘆旺店'
$ws.Range('H32').Value = '新北市蘆洲區'
$ws.Range('Q32').Value = 'THILF02958'
$ws.Range('R32').Value = '新北一'
$ws.Range('S32').Value = '吳宗鴻'
$ws.Range('T32').Value = 1
$ws.Range('U32').Value = '已完工'
$ws.Range('V32').Value = '2025-06-05 14:25:00'
$ws.Range('W32').Value = '2025-06-05 13:50:00'
$ws.Range('X32').Value = '2025-06-05 14:24:00'
$ws.Range('Z32').Value = 0.6
$ws.Range('AB32').Value = '到場處理'
$ws.Range('AC32').Value = 'PMQ2+EDC+STAR'
$ws.Range('AD32').Value = 'O'
$ws.Range('AJ32').Value = 'O'
$ws.Range('AK32').Value = 'O'

# ---- Row 33 ----
Apply-CenterStyle $ws.Range('A33:L33') $true
Apply-CenterStyle $ws.Range('N33') $true
Apply-CenterStyle $ws.Range('Q33:AB33') $true
Apply-CenterStyle $ws.Range('AD33:AK33') $true
Apply-LeftStyle $ws.Range('M33') $true $false
Apply-LeftStyle $ws.Range('O33') $true $false
Apply-LeftStyle $ws.Range('P33') $true $true
Apply-LeftStyle $ws.Range('AC33') $true $true
$ws.Range('A33').Value = 31
$ws.Range('B33').Value = '服務'
$ws.Range('C33').Value = 2025060799
$ws.Range('F33').Value = 4801
$ws.Range('G33').Value = '五股成州店'
$ws.Range('H33').Value = '新北市五股區'
$ws.Range('Q33').Value = 'THILF04801'
$ws.Range('R33').Value = '新北一'
$ws.Range('S33').Value = '湯家瑋'
$ws.Range('T33').Value = 1
$ws.Range('U33').Value = '已完工'
$ws.Range('V33').Value = '2025-06-05 14:36:37'
$ws.Range('W33').Value = '2025-06-05 14:00:00'
$ws.Range('X33').Value = '2025-06-05 14:40:00'
$ws.Range('Z33').Value = 0.7
$ws.Range('AB33').Value = '到場處理'
$ws.Range('AC33').Value = 'PMQ2+EDC+STAR'
$ws.Range('AD33').Value = 'O'
$ws.Range('AJ33').Value = 'O'
$ws.Range('AK33').Value = 'O'

# ---- Row 34 ----
Apply-CenterStyle $ws.Range('A34:L34') $false
Apply-CenterStyle $ws.Range('N34') $false
Apply-CenterStyle $ws.Range('Q34:AB34') $false
Apply-CenterStyle $ws.Range('AD34:AK34') $false
Apply-LeftStyle $ws.Range('M34') $false $false
Apply-LeftStyle $ws.Range('O34') $false $false
Apply-LeftStyle $ws.Range('P34') $false $true
Apply-LeftStyle $ws.Range('AC34') $false $true
$ws.Range('A34').Value = 32
$ws.Range('B34').Value = '服務'
$ws.Range('C34').Value = 2025060808
$ws.Range('F34').Value = 3416
$ws.Range('G34').Value = '板橋國慶店'
$ws.Range('H34').Value = '新北市板橋區'
$ws.Range('Q34').Value = 'THILF03416'
$ws.Range('R34').Value = '新北一'
$ws.Range('S34').Value = '狄澤洋'
$ws.Range('T34').Value = 1
$ws.Range('U34').Value = '已完工'
$ws.Range('V34').Value = '2025-06-05 15:00:25'
$ws.Range('W34').Value = '2025-06-05 14:30:00'
$ws.Range('X34').Value = '2025-06-05 14:55:00'
$ws.Range('Z34').Value = 0.4
$ws.Range('AB34').Value = '到場處理'
$ws.Range('AC34').Value = 'PMQ2+EDC+STAR'
$ws.Range('AD34').Value = 'O'
$ws.Range('AJ34').Value = 'O'
$ws.Range('AK34').Value = 'O'

# ---- Row 35 ----
Apply-CenterStyle $ws.Range('A35:L35') $true
Apply-CenterStyle $ws.Range('N35') $true
Apply-CenterStyle $ws.Range('Q35:AB35') $true
Apply-CenterStyle $ws.Range('AD35:AK35') $true
Apply-LeftStyle $ws.Range('M35') $true $false
Apply-LeftStyle $ws.Range('O35') $true $false
Apply-LeftStyle $ws.Range('P35') $true $true
Apply-LeftStyle $ws.Range('AC35') $true $true
$ws.Range('A35').Value = 33
$ws.Range('B35').Value = '服務'
$ws.Range('C35').Value = 2025060816
$ws.Range('F35').Value = 4609
$ws.Range('G35').Value = '蘆洲洲正店'
$ws.Range('H35').Value = '新北市蘆洲區'
$ws.Range('Q35').Value = 'THILF04609'
$ws.Range('R35').Value = '新北一'
$ws.Range('S35').Value = '吳宗鴻'
$ws.Range('T35').Value = 1
$ws.Range('U35').Value = '已完工'
$ws.Range('V35').Value = '2025-06-05 15:26:43'
$ws.Range('W35').Value = '2025-06-05 14:51:00'
$ws.Range('X35').Value = '2025-06-05 15:26:00'
$ws.Range('Z35').Value = 0.6
$ws.Range('AB35').Value = '到場處理'
$ws.Range('AC35').Value = 'PMQ2+EDC+STAR'
$ws.Range('AD35').Value = 'O'
$ws.Range('AJ35').Value = 'O'
$ws.Range('AK35').Value = 'O'

# ---- Row 36 ----
Apply-CenterStyle $ws.Range('A36:L36') $false
Apply-CenterStyle $ws.Range('N36') $false
Apply-CenterStyle $ws.Range('Q36:AB36') $false
Apply-CenterStyle $ws.Range('AD36:AK36') $false
Apply-LeftStyle $ws.Range('M36') $false $false
Apply-LeftStyle $ws.Range('O36') $false $false
Apply-LeftStyle $ws.Range('P36') $false $false
Apply-LeftStyle $ws.Range('AC36') $false $false
$ws.Range('A36').Value = 34
$ws.Range('B36').Value = '服務'
$ws.Range('C36').Value = 2025060817
$ws.Range('F36').Value = 4819
$ws.Range('G36').Value = '五股凌雲店'
$ws.Range('H36').Value = '新北市五股區'
$ws.Range('Q36').Value = 'THILF04819'
$ws.Range('R36').Value = '新北一'
$ws.Range('S36').Value = '湯家瑋'
$ws.Range('T36').Value = 1
$ws.Range('U36').Value = '已完工'
$ws.Range('V36').Value = '2025-06-05 15:27:15'
$ws.Range('W36').Value = '2025-06-05 14:50:00'
$ws.Range('X36').Value = '2025-06-05 15:30:00'
$ws.Range('Z36').Value = 0.7
$ws.Range('AB36').Value = '到場處理'
$ws.Range('AC36').Value = 'PMQ2+EDC+STAR'
$ws.Range('AD36').Value = 'O'
$ws.Range('AJ36').Value = 'O'
$ws.Range('AK36').Value = 'O'

# ---- Row 26 style tweak: P26, AC26 gain WrapText ----
Apply-LeftStyle $ws.Range('P26') $false $true
Apply-LeftStyle $ws.Range('AC26') $false $true

# ---- Print area, selection ----
$ws.PageSetup.PrintArea = '$A$1:$AK$36'
[void]$ws.Range('AC33').Select()
